$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the values of columns A, Q, R between row 4 and row 5.
$a4 = $ws.Range("A4").Value2
$q4 = $ws.Range("Q4").Value2
$r4 = $ws.Range("R4").Value2

$a5 = $ws.Range("A5").Value2
$q5 = $ws.Range("Q5").Value2
$r5 = $ws.Range("R5").Value2

$ws.Range("A4").Value = $a5
$ws.Range("Q4").Value = $q5
$ws.Range("R4").Value = $r5

$ws.Range("A5").Value = $a4
$ws.Range("Q5").Value = $q4
$ws.Range("R5").Value = $r4
